# Update the ObjTables schema header metadata (date + objTablesVersion)
# embedded in the first row(s) of every "!!<Class>" worksheet, mirroring
# a re-export of the SBtab workbook with a newer objtables library.

$wb = $excel.ActiveWorkbook

$oldDate = "2020-04-27 01:09:19"
$newDate = "2020-05-29 00:23:14"
$oldVersion = "0.0.9"
$newVersion = "1.0.0"

foreach ($ws in $wb.Worksheets) {
    # Row 1 holds the per-class "!!ObjTables ..." header on every sheet,
    # except the very first sheet, whose row 1 instead holds the
    # workbook-wide "!!!ObjTables ..." banner and whose per-class header
    # lives in row 2.
    for ($r = 1; $r -le 2; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value
        if ($val -ne $null -and $val -is [string] -and $val.StartsWith("!!")) {
            $newVal = $val.Replace($oldDate, $newDate).Replace("='" + $oldVersion + "'", "='" + $newVersion + "'")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
